# Apply updated PvsI model-fitting results to the respirometry rates sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New area/output-unit label used throughout (cm2 -> m2)
$newUnit = "umolO2/min/m2"

# Row-by-row target values: T (volume), V (area), Z (rate.abs),
# AB (rate.a.spec), AD (rate.output) -- AC (output.unit) is always $newUnit.
# $null marks a cell whose content must be cleared entirely (removed from the sheet).
$rowData = @{
    2  = @{ T = 0.1426048780487805;  V = 0.0002220611491829204; Z = -0.1272545377268952; AB = -573.0607906656857;  AD = -573.0607906656857 }
    3  = @{ T = 0.1423414634146342;  V = 0.0002565454225970831; Z = -0.1605930463497094; AB = -625.982895052189;   AD = -625.982895052189 }
    4  = @{ T = 0.1477268292682927;  V = 0.0002491214197856265; Z = -0.1284305000115185; AB = -515.533751060167;   AD = -515.533751060167 }
    5  = @{ T = 0.1470341463414634;  V = 0.0002082674398172554; Z = -0.1573829150353515; AB = -755.6770044009154;  AD = -755.6770044009154 }
    6  = @{ T = 0.1477560975609756;  V = 0.0001607142857142857; Z = -0.1931961599169131; AB = -1202.109439483015;  AD = -1202.109439483015 }
    7  = @{ T = 0.1455317073170732;  V = 0.0002247627833421192; Z = -0.106547518748761;  AB = -474.0443109150387;  AD = -474.0443109150387 }
    8  = @{ T = 0.1544;              V = 0;                     Z = 0;                   AB = $null;               AD = $null }
    9  = @{ T = 0.1426048780487805;  V = 0.0002220611491829204; Z = 0.1162725524360758;  AB = 523.6060106142095;   AD = 523.6060106142095 }
    10 = @{ T = 0.1423414634146342;  V = 0.0002565454225970831; Z = 0.09745889124492468; AB = 379.8894178594975;   AD = 379.8894178594975 }
    11 = @{ T = 0.1477268292682927;  V = 0.0002491214197856265; Z = 0.1748284922790208;  AB = 701.7802500863388;   AD = 701.7802500863388 }
    12 = @{ T = 0.1470341463414634;  V = 0.0002082674398172554; Z = 0.2180113928172306;  AB = 1046.785772219244;   AD = 1046.785772219244 }
    13 = @{ T = 0.1477560975609756;  V = 0.0001607142857142857; Z = 0.1851069083597895;  AB = 1151.776318683135;   AD = 1151.776318683135 }
    14 = @{ T = 0.1455317073170732;  V = 0.0002247627833421192; Z = 0.1029625439693433;  AB = 458.094273608547;    AD = 458.094273608547 }
    15 = @{ T = 0.1544;              V = 0;                     Z = 0.00008312341187081832; AB = "Inf";            AD = "Inf" }
}

foreach ($r in $rowData.Keys) {
    $d = $rowData[$r]

    $ws.Range("T$r").Value = $d.T
    $ws.Range("V$r").Value = $d.V
    $ws.Range("Z$r").Value = $d.Z

    if ($null -eq $d.AB) {
        $ws.Range("AB$r").ClearContents()
    } else {
        $ws.Range("AB$r").Value = $d.AB
    }

    $ws.Range("AC$r").Value = $newUnit

    if ($null -eq $d.AD) {
        $ws.Range("AD$r").ClearContents()
    } else {
        $ws.Range("AD$r").Value = $d.AD
    }
}
